# SQL DATABASE AUTOMATION TESTING
# Update the sample employee rows on the "addEmployeeData" sheet with a
# fresh batch of test data (rows 2-4), then leave the selection on E4
# the way the author left it after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("addEmployeeData")

# Row 2 keeps its name/username, only the password test value changes.
$ws.Range("E2").Value = "JohKhan_*@123"

# Row 3 gets a new set of test values.
$ws.Range("A3").Value = "Haryy"
$ws.Range("B3").Value = "Mimis"
$ws.Range("C3").Value = "Annys"
$ws.Range("D3").Value = "mary567298045"
$ws.Range("E3").Value = "mayrKhan_*@123"

# Row 4 gets a new set of test values.
$ws.Range("A4").Value = "Jordon"
$ws.Range("B4").Value = "Jelli"
$ws.Range("C4").Value = "John"
$ws.Range("D4").Value = "Jor35323765"
$ws.Range("E4").Value = "jorKhan_*@123"

# Match the author's final cursor position on the sheet.
$ws.Activate() | Out-Null
$ws.Range("E4").Select() | Out-Null
